$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns remain stored as text,
# matching the workbook's original inline-string cell type, rather than
# letting Excel auto-convert numeric-looking / percent-looking text.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("E2").Value = "-3.55%"
$ws.Range("D3").Value = "37.43"
$ws.Range("E3").Value = "-5.73%"
$ws.Range("D4").Value = "5.087"
$ws.Range("E4").Value = "-0.96%"
$ws.Range("D5").Value = "0.07728"
$ws.Range("E6").Value = "0.65%"
$ws.Range("E7").Value = "-8.40%"
$ws.Range("D8").Value = "8.204"
$ws.Range("E8").Value = "-1.73%"
$ws.Range("D9").Value = "3.018"
$ws.Range("E9").Value = "-6.62%"
$ws.Range("D10").Value = "0.9183"
$ws.Range("E10").Value = "-2.34%"
$ws.Range("D11").Value = "0.1150"
$ws.Range("E11").Value = "-15.50%"
$ws.Range("D12").Value = "0.1883"
$ws.Range("E12").Value = "-5.33%"
$ws.Range("D13").Value = "0.08732"
$ws.Range("E13").Value = "-4.44%"
$ws.Range("D14").Value = "0.03404"
$ws.Range("E14").Value = "-2.65%"
$ws.Range("D15").Value = "0.09707"
$ws.Range("E15").Value = "-0.83%"
$ws.Range("D16").Value = "0.001365"
$ws.Range("E16").Value = "-2.79%"
$ws.Range("D17").Value = "0.005975"
$ws.Range("E17").Value = "-0.61%"
$ws.Range("D18").Value = "3.585"
$ws.Range("E18").Value = "-2.75%"
$ws.Range("D19").Value = "0.3407"
$ws.Range("E20").Value = "-2.92%"
$ws.Range("D21").Value = "5.021"
$ws.Range("E21").Value = "0.85%"
$ws.Range("E22").Value = "5.99%"
$ws.Range("E23").Value = "5,171.33%"
$ws.Range("D24").Value = "0.04325"
$ws.Range("E25").Value = "-1.29%"
$ws.Range("D26").Value = "0.004540"
$ws.Range("E26").Value = "-5.56%"
$ws.Range("E27").Value = "3.89%"
$ws.Range("D39").Value = "0.02207"
$ws.Range("E39").Value = "-2.20%"
$ws.Range("D40").Value = "0.04918"
$ws.Range("E40").Value = "-5.70%"
$ws.Range("D41").Value = "0.007562"
$ws.Range("E41").Value = "-2.51%"
$ws.Range("D42").Value = "0.009944"
$ws.Range("E42").Value = "0.88%"
$ws.Range("D43").Value = "0.1335"
$ws.Range("E43").Value = "-5.07%"
$ws.Range("D44").Value = "0.001997"
$ws.Range("E44").Value = "-2.39%"
$ws.Range("D45").Value = "0.008782"
$ws.Range("E45").Value = "-5.89%"
$ws.Range("D46").Value = "0.00006540"
$ws.Range("E46").Value = "-1.04%"
$ws.Range("E47").Value = "0.23%"
$ws.Range("D48").Value = "0.003000"
$ws.Range("E48").Value = "1.87%"
$ws.Range("E49").Value = "-22.92%"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").Value = "0.23%"
$ws.Range("E51").Value = "0.23%"
